# Sprint5_Backlog.xlsx - "Add files via upload"
# Adds two new columns (E: "Reviewers 1", F: "Reviewers 2") with reviewer
# names for each user-story row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row
$ws.Range("E1").Value = "Reviewers 1"
$ws.Range("F1").Value = "Reviewers 2"

# Row -> (Reviewers 1, Reviewers 2)
$reviewers = @{
    2  = @("mohamed atef",  "omar fatayry")
    3  = @("mohamed atef",  "omar fatayry")
    4  = @("khaled moataz", "youssef abo bakr")
    5  = @("khaled moataz", "youssef abo bakr")
    6  = @("mohamed alaa",  "andrew")
    7  = @("mohamed alaa",  "andrew")
    8  = @("waly",          "ahmed hesham")
    9  = @("waly",          "ahmed hesham")
    10 = @("mohamed atef",  "omar fatayry")
    11 = @("bassem",        "reyad")
    12 = @("bassem",        "reyad")
    13 = @("bassem",        "reyad")
    14 = @("khaled moataz", "youssef abo bakr")
    16 = @("mohamed alaa",  "andrew")
    17 = @("waly",          "ahmed hesham")
    18 = @("khaled moataz", "youssef abo bakr")
    19 = @("waly",          $null)
}

foreach ($row in $reviewers.Keys) {
    $pair = $reviewers[$row]
    $ws.Cells.Item($row, 5).Value = $pair[0]
    if ($pair[1]) {
        $ws.Cells.Item($row, 6).Value = $pair[1]
    }
}

# Reflect the updated selection (matches the saved view state in the workbook)
$ws.Range("F9").Select()
